$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove stray leftover formulas in row 92 (E92 / I92), leaving the cells blank ---
$ws.Range("E92").ClearContents()
$ws.Range("I92").ClearContents()

# --- Insert a new table row after row 95 (i.e. at row 96), shifting all following rows down ---
$ws.Rows.Item(96).Insert()

# Grow Table1 to include the newly inserted row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K148"))

# Copy formatting from row 95 into the freshly inserted row 96
$ws.Range("A95:K95").Copy()
$ws.Range("A96:K96").PasteSpecial(-4122)

# Restore the calculated column formula in the EARNED helper column (G) for row 96
$ws.Range("G96").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Row 95 (February 2023 period): SL(2-0-0) leave, 2 hrs absence w/pay, remarks date
$ws.Range("B95").Value = "SL(2-0-0)"
$ws.Range("H95").Value = 2
$ws.Range("K95").Value = "2/6,9/2023"

# Row 96 (new continuation row, still within February 2023 period): SL(2-0-0) leave, 2 hrs absence w/pay, remarks date
$ws.Range("B96").Value = "SL(2-0-0)"
$ws.Range("H96").Value = 2
$ws.Range("K96").Value = "2/13,14/2023"

# Row 97 (March 2023 period, shifted down from the former row 96): EARNED value recorded
$ws.Range("C97").Value = 1.25
